# Auto-generated Excel COM-interop script applying the Sagittarius_Profits.xlsx diff
# Updates scheduled-runner price/profit values across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 819.13635
$ws.Range("I2").Value = 214.5625
$ws.Range("J2").Value = 2431.3333
$ws.Range("K2").Value = 214.5625
$ws.Range("L2").Value = 2431.3333
$ws.Range("M2").Value = -101.5625
$ws.Range("N2").Value = -2657.3333

$ws.Range("H28").Value = 2480.8333
$ws.Range("J28").Value = 1150
$ws.Range("L28").Value = 1150
$ws.Range("N28").Value = -2120

$ws.Range("H33").Value = 226.4
$ws.Range("I33").Value = 153.66667
$ws.Range("K33").Value = 153.66667
$ws.Range("M33").Value = 75.33332999999999

$ws.Range("H74").Value = 74353.125
$ws.Range("I74").Value = 77042.39
$ws.Range("K74").Value = 77042.39
$ws.Range("M74").Value = -76106.39

$ws.Range("H77").Value = 74353.125
$ws.Range("I77").Value = 77042.39
$ws.Range("K77").Value = 385211.95
$ws.Range("M77").Value = -380531.95

$ws.Range("H100").Value = 2627.5715
$ws.Range("I100").Value = 2220.6
$ws.Range("J100").Value = 3645
$ws.Range("K100").Value = 2220.6
$ws.Range("L100").Value = 3645
$ws.Range("M100").Value = -1679.6
$ws.Range("N100").Value = -4727

$ws.Range("H112").Value = 3516.1428
$ws.Range("J112").Value = 3502.1667
$ws.Range("L112").Value = 10506.5001
$ws.Range("N112").Value = -12722.5001

$ws.Range("H132").Value = 1752.5
$ws.Range("I132").Value = 1752.5
$ws.Range("K132").Value = 5257.5
$ws.Range("M132").Value = -2727.5

$ws.Range("H137").Value = 1374.5
$ws.Range("J137").Value = 1249.6666
$ws.Range("L137").Value = 3748.9998
$ws.Range("N137").Value = -8848.9998

$ws.Range("H138").Value = 2531.8096
$ws.Range("I138").Value = 2893.3333
$ws.Range("J138").Value = 2170.2856
$ws.Range("K138").Value = 8679.999899999999
$ws.Range("L138").Value = 6510.8568
$ws.Range("M138").Value = -3539.999899999999
$ws.Range("N138").Value = -16790.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5111.684
$ws.Range("I32").Value = 4840.1665
$ws.Range("K32").Value = 4840.1665
$ws.Range("M32").Value = -4553.1665

$ws.Range("H45").Value = 5826.375
$ws.Range("I45").Value = 6601.8335
$ws.Range("K45").Value = 6601.8335
$ws.Range("M45").Value = -6224.8335

$ws.Range("H74").Value = 1766.3334
$ws.Range("I74").Value = 1399.5
$ws.Range("K74").Value = 1399.5
$ws.Range("M74").Value = -525.5

$ws.Range("H77").Value = 1766.3334
$ws.Range("I77").Value = 1399.5
$ws.Range("K77").Value = 6997.5
$ws.Range("M77").Value = -2629.5

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H122").Value = 1549.8889
$ws.Range("I122").Value = 924.8333
$ws.Range("K122").Value = 2774.4999
$ws.Range("M122").Value = -324.4998999999998

$ws.Range("H132").Value = 1708.875
$ws.Range("I132").Value = 1708.875
$ws.Range("K132").Value = 5126.625
$ws.Range("M132").Value = -2596.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1335.5
$ws.Range("I20").Value = 1199
$ws.Range("J20").Value = 1381
$ws.Range("K20").Value = 1199
$ws.Range("L20").Value = 1381
$ws.Range("M20").Value = -952
$ws.Range("N20").Value = -1875

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").ClearContents()

$ws.Range("H82").Value = 8583
$ws.Range("I82").Value = 8583
$ws.Range("K82").Value = 8583
$ws.Range("M82").Value = -8200

$ws.Range("H85").Value = 8583
$ws.Range("I85").Value = 8583
$ws.Range("K85").Value = 8583
$ws.Range("M85").Value = -7257

$ws.Range("H105").Value = 2019.409
$ws.Range("I105").Value = 1864.6316
$ws.Range("J105").Value = 2999.6667
$ws.Range("K105").Value = 1864.6316
$ws.Range("L105").Value = 2999.6667
$ws.Range("M105").Value = -117.6315999999999
$ws.Range("N105").Value = -6493.6667

$ws.Range("H134").Value = 1662.8334
$ws.Range("I134").Value = 1662.8334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4988.5002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2453.5002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4077.25
$ws.Range("I16").Value = 4437
$ws.Range("J16").Value = 2998
$ws.Range("K16").Value = 4437
$ws.Range("L16").Value = 2998
$ws.Range("M16").Value = -4150
$ws.Range("N16").Value = -3572

$ws.Range("H58").Value = 2748.5386
$ws.Range("I58").Value = 2744.25
$ws.Range("K58").Value = 2744.25
$ws.Range("M58").Value = -2541.25

$ws.Range("H86").Value = 9351.429
$ws.Range("I86").Value = 8863.5
$ws.Range("K86").Value = 8863.5
$ws.Range("M86").Value = -7740.5

$ws.Range("H89").Value = 9351.429
$ws.Range("I89").Value = 8863.5
$ws.Range("K89").Value = 44317.5
$ws.Range("M89").Value = -38701.5

$ws.Range("H105").Value = 3080.6155
$ws.Range("I105").Value = 2316
$ws.Range("J105").Value = 3845.2307
$ws.Range("K105").Value = 2316
$ws.Range("L105").Value = 3845.2307
$ws.Range("M105").Value = -569
$ws.Range("N105").Value = -7339.2307

$ws.Range("H113").Value = 4077.25
$ws.Range("I113").Value = 4437
$ws.Range("J113").Value = 2998
$ws.Range("K113").Value = 4437
$ws.Range("L113").Value = 2998
$ws.Range("M113").Value = -2267
$ws.Range("N113").Value = -7338

$ws.Range("H122").Value = 3260
$ws.Range("I122").Value = 3260
$ws.Range("K122").Value = 9780
$ws.Range("M122").Value = -7330

$ws.Range("H134").Value = 1076.4445
$ws.Range("I134").Value = 914
$ws.Range("K134").Value = 2742
$ws.Range("M134").Value = -207

$ws.Range("H136").Value = 2748.5386
$ws.Range("I136").Value = 2744.25
$ws.Range("K136").Value = 8232.75
$ws.Range("M136").Value = -5682.75

$ws.Range("H141").Value = 641666.7
$ws.Range("J141").Value = 641666.7
$ws.Range("L141").Value = 641666.7
$ws.Range("N141").Value = -652026.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 14083.333
$ws.Range("J57").Value = 16300
$ws.Range("L57").Value = 48900
$ws.Range("N57").Value = -50018

$ws.Range("H126").Value = 3030
$ws.Range("I126").Value = 3030
$ws.Range("K126").Value = 9090
$ws.Range("M126").Value = -4150

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 26670672
$ws.Range("I36").Value = 6009
$ws.Range("J36").Value = 80000000
$ws.Range("K36").Value = 6009
$ws.Range("L36").Value = 80000000
$ws.Range("M36").Value = -5524
$ws.Range("N36").Value = -80000970

$ws.Range("H70").Value = 7049.3447
$ws.Range("I70").Value = 6801.0557
$ws.Range("J70").Value = 7455.636
$ws.Range("K70").Value = 6801.0557
$ws.Range("L70").Value = 7455.636
$ws.Range("M70").Value = -6531.0557
$ws.Range("N70").Value = -7995.636

$ws.Range("H73").Value = 7049.3447
$ws.Range("I73").Value = 6801.0557
$ws.Range("J73").Value = 7455.636
$ws.Range("K73").Value = 6801.0557
$ws.Range("L73").Value = 7455.636
$ws.Range("M73").Value = -5865.0557
$ws.Range("N73").Value = -9327.636

$ws.Range("H102").Value = 1749.2222
$ws.Range("I102").Value = 1749.2222
$ws.Range("K102").Value = 1749.2222
$ws.Range("M102").Value = -127.2221999999999

$ws.Range("H132").Value = 789.8
$ws.Range("I132").Value = 789.8
$ws.Range("K132").Value = 2369.4
$ws.Range("M132").Value = 160.6000000000004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2794.2273
$ws.Range("I40").Value = 2242.5715
$ws.Range("K40").Value = 2242.5715
$ws.Range("M40").Value = -2106.5715

$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()

$ws.Range("H132").Value = 2632.25
$ws.Range("I132").Value = 2416.5454
$ws.Range("K132").Value = 7249.6362
$ws.Range("M132").Value = -4719.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1253287.5
$ws.Range("I81").Value = 4559.6
$ws.Range("K81").Value = 9119.200000000001
$ws.Range("M81").Value = -8058.200000000001

$ws.Range("H84").Value = 1253287.5
$ws.Range("I84").Value = 4559.6
$ws.Range("K84").Value = 45596
$ws.Range("M84").Value = -40292

$ws.Range("H107").Value = 1242.1111
$ws.Range("I107").Value = 1254.5714
$ws.Range("K107").Value = 3763.7142
$ws.Range("M107").Value = -1843.7142

$ws.Range("H132").Value = 3855
$ws.Range("I132").Value = 4217.2856
$ws.Range("K132").Value = 12651.8568
$ws.Range("M132").Value = -10121.8568
